# Natmi following Dr Hou advice
# Update Rspo3-Lgr5 LR-pair sheet: refresh stats for existing rows (FAPs->FAPs, FAPs->sCs)
# and add two new target-cluster rows (FAPs->Neutro, FAPs->sCs duplicate) plus rename
# the former sCs row's target cluster to M1, per the updated NATMI run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 : Sending=FAPs, Ligand=Rspo3, Receptor=Lgr5, Target=FAPs ---
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.822099333333334
$ws.Range("H2").Value = 8.466298
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.6233136666666667
$ws.Range("N2").Value = 1.869941
$ws.Range("O2").Value = 0.1671957790227948
$ws.Range("P2").Value = 0.1671957790227948
$ws.Range("Q2").Value = 1.759053083157556
$ws.Range("R2").Value = 15.831477748418
$ws.Range("S2").Value = 0.1671957790227948
$ws.Range("T2").Value = 0.1671957790227948

# --- Row 3 : Sending=FAPs, Ligand=Rspo3, Receptor=Lgr5, Target=M1 (was sCs) ---
$ws.Range("D3").Value = "M1"
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.822099333333334
$ws.Range("H3").Value = 8.466298
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.190009
$ws.Range("N3").Value = 0.570027
$ws.Range("O3").Value = 0.05096744139468926
$ws.Range("P3").Value = 0.05096744139468926
$ws.Range("Q3").Value = 0.5362242722273334
$ws.Range("R3").Value = 4.826018450046
$ws.Range("S3").Value = 0.05096744139468926
$ws.Range("T3").Value = 0.05096744139468926

# --- Row 4 (new) : Sending=FAPs, Ligand=Rspo3, Receptor=Lgr5, Target=Neutro ---
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Rspo3"
$ws.Range("C4").Value = "Lgr5"
$ws.Range("D4").Value = "Neutro"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.822099333333334
$ws.Range("H4").Value = 8.466298
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.022719
$ws.Range("N4").Value = 0.068157
$ws.Range("O4").Value = 0.006094076075585607
$ws.Range("P4").Value = 0.006094076075585606
$ws.Range("Q4").Value = 0.064115274754
$ws.Range("R4").Value = 0.5770374727859999
$ws.Range("S4").Value = 0.006094076075585607
$ws.Range("T4").Value = 0.006094076075585606

# --- Row 5 (new) : Sending=FAPs, Ligand=Rspo3, Receptor=Lgr5, Target=sCs ---
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Rspo3"
$ws.Range("C5").Value = "Lgr5"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.822099333333334
$ws.Range("H5").Value = 8.466298
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.892005
$ws.Range("N5").Value = 8.676015
$ws.Range("O5").Value = 0.7757427035069303
$ws.Range("P5").Value = 0.7757427035069303
$ws.Range("Q5").Value = 8.161525382496666
$ws.Range("R5").Value = 73.45372844246999
$ws.Range("S5").Value = 0.7757427035069303
$ws.Range("T5").Value = 0.7757427035069303
